$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Col8a1", "Itga2", "ECs", 3, 1, 4.656903666666667, 13.970711, 0.03472688986918189, 0.03526504803992261, 3, 1, 3.042935333333334, 9.128806, 0.5356624298519768, 0.5831288282815508, 14.17065671122956, 127.535910401066, 0.01860189020852797, 0.02056406614281267),
  @("ECs", "Col8a1", "Itga2", "FAPs", 3, 1, 4.656903666666667, 13.970711, 0.03472688986918189, 0.03526504803992261, 3, 1, 1.067915666666667, 3.203747, 0.1879902916822836, 0.2046485853922773, 4.973180383790778, 44.758623454117, 0.006528318155726044, 0.007216942195160863),
  @("ECs", "Col8a1", "Itga2", "M1", 3, 1, 4.656903666666667, 13.970711, 0.03472688986918189, 0.03526504803992261, 1, 0.3333333333333333, 0.02258333333333333, 0.06775, 0.003975451950942043, 0.004327726771285868, 0.1051684078055556, 0.94651567025, 0.0001380550820805886, 0.0001526174924930553),
  @("ECs", "Col8a1", "Itga2", "M2", 3, 1, 4.656903666666667, 13.970711, 0.03472688986918189, 0.03526504803992261, 3, 1, 0.160044, 0.480132, 0.02817330916767092, 0.03066981712400039, 0.745309490428, 6.707785413852, 0.0009783714047161204, 0.001081572574253515),
  @("ECs", "Col8a1", "Itga2", "sCs", 3, 1, 4.656903666666667, 13.970711, 0.03472688986918189, 0.03526504803992261, 2, 1, 1.3872175, 2.774435, 0.2441985173471266, 0.1772250424308857, 6.460138262214167, 38.760829573285, 0.00848025501813117, 0.006249849635202507),
  @("FAPs", "Col8a1", "Itga2", "ECs", 3, 1, 123.304606, 369.913818, 0.9194919584818978, 0.9337412077596615, 3, 1, 3.042935333333334, 9.128806, 0.5356624298519768, 0.5831288282815508, 375.2079423601454, 3376.871481241308, 0.4925372967097663, 0.5444914163990915),
  @("FAPs", "Col8a1", "Itga2", "FAPs", 3, 1, 123.304606, 369.913818, 0.9194919584818978, 0.9337412077596615, 3, 1, 1.067915666666667, 3.203747, 0.1879902916822836, 0.2046485853922773, 131.6789205195607, 1185.110284676046, 0.1728555614745262, 0.1910888172904912),
  @("FAPs", "Col8a1", "Itga2", "M1", 3, 1, 123.304606, 369.913818, 0.9194919584818978, 0.9337412077596615, 1, 0.3333333333333333, 0.02258333333333333, 0.06775, 0.003975451950942043, 0.004327726771285868, 2.784629018833333, 25.0616611695, 0.00365539610022238, 0.004040976822274287),
  @("FAPs", "Col8a1", "Itga2", "M2", 3, 1, 123.304606, 369.913818, 0.9194919584818978, 0.9337412077596615, 3, 1, 0.160044, 0.480132, 0.02817330916767092, 0.03066981712400039, 19.734162362664, 177.607461263976, 0.02590513122349774, 0.02863767208313207),
  @("FAPs", "Col8a1", "Itga2", "sCs", 3, 1, 123.304606, 369.913818, 0.9194919584818978, 0.9337412077596615, 2, 1, 1.3872175, 2.774435, 0.2441985173471266, 0.1772250424308857, 171.050307273805, 1026.30184364283, 0.2245385729738851, 0.1654823251646725),
  @("sCs", "Col8a1", "Itga2", "ECs", 2, 1, 6.139289, 12.278578, 0.04578115164892033, 0.03099374420041592, 3, 1, 3.042935333333334, 9.128806, 0.5356624298519768, 0.5831288282815508, 18.68145941964467, 112.088756517868, 0.0245232429336825, 0.01807334573964664),
  @("sCs", "Col8a1", "Itga2", "FAPs", 2, 1, 6.139289, 12.278578, 0.04578115164892033, 0.03099374420041592, 3, 1, 1.067915666666667, 3.203747, 0.1879902916822836, 0.2046485853922773, 6.556242905294333, 39.337457431766, 0.008606412052031394, 0.006342825906625216),
  @("sCs", "Col8a1", "Itga2", "M1", 2, 1, 6.139289, 12.278578, 0.04578115164892033, 0.03099374420041592, 1, 0.3333333333333333, 0.02258333333333333, 0.06775, 0.003975451950942043, 0.004327726771285868, 0.1386456099166667, 0.8318736595, 0.0001820007686390739, 0.0001341324565185261),
  @("sCs", "Col8a1", "Itga2", "M2", 2, 1, 6.139289, 12.278578, 0.04578115164892033, 0.03099374420041592, 3, 1, 0.160044, 0.480132, 0.02817330916767092, 0.03066981712400039, 0.9825563687159999, 5.895338212296, 0.00128980653945706, 0.0009505724666148038),
  @("sCs", "Col8a1", "Itga2", "sCs", 2, 1, 6.139289, 12.278578, 0.04578115164892033, 0.03099374420041592, 2, 1, 1.3872175, 2.774435, 0.2441985173471266, 0.1772250424308857, 8.5165291383575, 34.06611655343, 0.01117968935511031, 0.005492867631010729)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

Write-Output "Applied $($data.Count) rows"